# The two pairs of observation rows (2<->4 and 3<->5) had their Id,
# Antal (count), Ost (easting) and Nord (northing) values swapped between
# each other. Swap them back into place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-NumericCells {
    param($ws, $addr1, $addr2)
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

function Swap-TextCells {
    param($ws, $addr1, $addr2, $scratch)
    $ws.Range($addr1).Copy() | Out-Null
    $ws.Range($scratch).PasteSpecial(-4163) | Out-Null
    $ws.Range($addr2).Copy() | Out-Null
    $ws.Range($addr1).PasteSpecial(-4163) | Out-Null
    $ws.Range($scratch).Copy() | Out-Null
    $ws.Range($addr2).PasteSpecial(-4163) | Out-Null
}

# Row 2 <-> Row 4: Id (A), Ost (Q), Nord (R)
Swap-NumericCells $ws "A2" "A4"
Swap-NumericCells $ws "Q2" "Q4"
Swap-NumericCells $ws "R2" "R4"

# Row 3 <-> Row 5: Id (A), Antal (I, text-typed), Ost (Q), Nord (R)
Swap-NumericCells $ws "A3" "A5"
Swap-TextCells $ws "I3" "I5" "AZ1"
Swap-NumericCells $ws "Q3" "Q5"
Swap-NumericCells $ws "R3" "R5"

$ws.Range("AZ1").Clear() | Out-Null
$excel.CutCopyMode = 0
